$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.878652691841125
$ws.Range("B1").Value = 3.532402992248535
$ws.Range("C1").Value = 3.150560617446899
$ws.Range("D1").Value = 2.573280572891235
$ws.Range("E1").Value = 1.655121564865112
